$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '93.590.97'
$ws.Range('E2').Value = '  -5.22%  '
$ws.Range('D3').Value = '3.386.63'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.78'
$ws.Range('E5').Value = '  -8.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '631.47'
$ws.Range('E6').Value = '  -6.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.40'
$ws.Range('E7').Value = '  -6.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.393'
$ws.Range('E8').Value = '  -9.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  +0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.945'
$ws.Range('E10').Value = '  -11.02%  '
$ws.Range('D11').Value = '3.385.67'
$ws.Range('E11').Value = '  -1.35%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.196'
$ws.Range('E12').Value = '  -6.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.64'
$ws.Range('E13').Value = '  -13.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.06'
$ws.Range('E14').Value = '  -2.37%  '
$ws.Range('D15').Value = '93.372.28'
$ws.Range('E15').Value = '  -5.22%  '
$ws.Range('D16').Value = '4.018.03'
$ws.Range('E16').Value = '  -1.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000245'
$ws.Range('E17').Value = '  -6.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.06'
$ws.Range('E18').Value = '  -11.71%  '
$ws.Range('D19').Value = '3.380.20'
$ws.Range('E19').Value = '  -1.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.05'
$ws.Range('E20').Value = '  -8.63%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.00'
$ws.Range('E21').Value = '  -5.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '492.20'
$ws.Range('E22').Value = '  -6.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.457'
$ws.Range('E23').Value = '  -14.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.17'
$ws.Range('E24').Value = '  -8.42%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000187'
$ws.Range('E25').Value = '  -8.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.34'
$ws.Range('E26').Value = '  -8.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '91.97'
$ws.Range('E27').Value = '  -6.54%  '
$ws.Range('B28').Value = 'Aptos'
$ws.Range('C28').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.60'
$ws.Range('E28').Value = '  -8.98%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '11.39'
$ws.Range('E29').Value = '  -7.42%  '
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.67'
$ws.Range('E31').Value = '  -10.54%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.132'
$ws.Range('E32').Value = '  -9.69%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.979'
$ws.Range('E33').Value = '  -2.30%  '
$ws.Range('B34').Value = 'Cronos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.172'
$ws.Range('E34').Value = '  -9.55%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '28.81'
$ws.Range('E35').Value = '  -2.83%  '
$ws.Range('B36').Value = 'PolygonEcosystemToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.534'
$ws.Range('E36').Value = '  -7.42%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '534.02'
$ws.Range('E37').Value = '  -0.90%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.51'
$ws.Range('E38').Value = '  -7.84%  '
$ws.Range('B39').Value = 'Fetch.AI'
$ws.Range('C39').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.41'
$ws.Range('E39').Value = '  -7.83%  '
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.149'
$ws.Range('E41').Value = '  -5.10%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.886'
$ws.Range('E42').Value = '  +0.77%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '24.04'
$ws.Range('E43').Value = '  -1.65%  '
$ws.Range('B44').Value = 'MantraDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.65'
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.67'
$ws.Range('E45').Value = '  -7.76%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.52'
$ws.Range('E46').Value = '  -3.78%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.16'
$ws.Range('E47').Value = '  -4.36%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0397'
$ws.Range('E48').Value = '  -8.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '52.99'
$ws.Range('E49').Value = '  -5.14%  '
$ws.Range('B50').Value = 'dogwifhat'
$ws.Range('C50').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.15'
$ws.Range('E50').Value = '  -3.80%  '
$ws.Range('B51').Value = 'Cosmos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.96'
$ws.Range('E51').Value = '  -9.20%  '
